$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.742.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "'2.268.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'248.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Value = "'0.642"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'77.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.69%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.651"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'39.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "'0.0969"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "'7.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "'2.608.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "'0.866"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "'2.261.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "'42.638.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "'0.0₃0990"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").Value = "'71.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'232.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "'2.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("E24").Value = "  -5.72%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'11.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("D27").Value = "'2.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.84%  "
$ws.Range("D28").Value = "'2.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "'167.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("E32").Value = "  +5.81%  "
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").Value = "'30.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.66%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "'0.0303"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.79%  "
$ws.Range("D39").Value = "'13.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("D41").Value = "'5.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "'0.207"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'110.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.42%  "
$ws.Range("D44").Value = "'60.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("D45").Value = "'8.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "'4.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.99%  "
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "'4.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.66%  "
